# Auto update Excel log
# Appends 6 new PRESENCE_DETECTED log rows (rows 40-45) to the "mmWave" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("mmWave")

$newRows = @(
    @("2026-02-01", "17:39:50", "17:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-02-01", "17:40:17", "17:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-02-01", "17:40:28", "17:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-02-01", "17:40:38", "17:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-02-01", "17:40:49", "17:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-02-01", "17:40:59", "17:00", "Living Room", "PRESENCE_DETECTED", "Active")
)

$startRow = 40
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]

    # Column A holds a date-shaped string ("2026-02-01"); a leading apostrophe
    # forces Excel to keep it as literal text instead of auto-converting it
    # to a date serial number, matching the rest of the log.
    $ws.Cells.Item($r, 1).Value = "'" + $rowData[0]
    $ws.Cells.Item($r, 2).Value = $rowData[1]
    $ws.Cells.Item($r, 3).Value = $rowData[2]
    $ws.Cells.Item($r, 4).Value = $rowData[3]
    $ws.Cells.Item($r, 5).Value = $rowData[4]
    $ws.Cells.Item($r, 6).Value = $rowData[5]

    # Clear the quote-prefix style side effect so the cell ends up with the
    # plain/default style, just like all the other log rows.
    $ws.Range("A" + $r + ":F" + $r).Style = "Normal"
}
